$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 281, shifting existing rows 281..344 down to 282..345
$ws.Rows.Item(281).Insert()

# Populate the newly inserted row 281 with the new record
$ws.Range("A281").Value = 4
$ws.Range("B281").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C281").Value = "Los Lagos"
$ws.Range("D281").Value = 45275
$ws.Range("E281").Value = 10
$ws.Range("F281").Value = "Fruta"
$ws.Range("G281").Value = 100103
$ws.Range("H281").Value = "Frutos de hueso (carozo)"
$ws.Range("I281").Value = 100103002
$ws.Range("J281").Value = "Ciruela"
$ws.Range("K281").Value = "Black Amber"
$ws.Range("L281").Value = "Primera"
$ws.Range("M281").Value = 350
$ws.Range("N281").Value = 18000
$ws.Range("O281").Value = 18000
$ws.Range("P281").Value = 18000
$ws.Range("Q281").Value = "$/caja 14 kilos granel"
$ws.Range("R281").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S281").Value = 1286
$ws.Range("T281").Value = 14
